$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply revised localization/comparison results (nfeature:512, hamming_tolerance:50, k_knn:15)
# Updates existing rows with recalculated metrics, and appends a new row 45 (TD_4.png).

$ws.Cells.Item(1, 2).Value = 0.462
$ws.Cells.Item(1, 3).Value = 0.002
$ws.Cells.Item(2, 2).Value = 0.864
$ws.Cells.Item(3, 2).Value = 0.846
$ws.Cells.Item(3, 3).Value = 0.003
$ws.Cells.Item(4, 2).Value = 0.451
$ws.Cells.Item(5, 2).Value = 0.452
$ws.Cells.Item(5, 3).Value = 0.002
$ws.Cells.Item(6, 2).Value = 1.228
$ws.Cells.Item(6, 3).Value = 0.004
$ws.Cells.Item(6, 4).Value = 0.599
$ws.Cells.Item(6, 5).Value = 0.9330000000000001
$ws.Cells.Item(7, 2).Value = 1.08
$ws.Cells.Item(7, 3).Value = 0.004
$ws.Cells.Item(7, 4).Value = 0.751
$ws.Cells.Item(7, 5).Value = 0.9330000000000001
$ws.Cells.Item(8, 2).Value = 1.327
$ws.Cells.Item(9, 2).Value = 1.096
$ws.Cells.Item(9, 3).Value = 0.004
$ws.Cells.Item(9, 4).Value = 0.73
$ws.Cells.Item(9, 5).Value = 0.867
$ws.Cells.Item(10, 2).Value = 1.232
$ws.Cells.Item(10, 3).Value = 0.004
$ws.Cells.Item(10, 4).Value = 0.734
$ws.Cells.Item(10, 5).Value = 0.867
$ws.Cells.Item(10, 6).Value = 'Muhammad Iqbal Baqi'
$ws.Cells.Item(10, 7).Value = 'Benar'
$ws.Cells.Item(11, 2).Value = 0.831
$ws.Cells.Item(11, 3).Value = 0.003
$ws.Cells.Item(12, 2).Value = 1.095
$ws.Cells.Item(12, 3).Value = 0.004
$ws.Cells.Item(13, 2).Value = 0.834
$ws.Cells.Item(13, 3).Value = 0.003
$ws.Cells.Item(13, 4).Value = 0.791
$ws.Cells.Item(14, 2).Value = 0.922
$ws.Cells.Item(15, 2).Value = 0.952
$ws.Cells.Item(15, 4).Value = 0.838
$ws.Cells.Item(15, 5).Value = 0.667
$ws.Cells.Item(16, 2).Value = 0.761
$ws.Cells.Item(16, 3).Value = 0.003
$ws.Cells.Item(17, 2).Value = 0.702
$ws.Cells.Item(17, 4).Value = 0.822
$ws.Cells.Item(17, 5).Value = 0.8
$ws.Cells.Item(18, 2).Value = 1.025
$ws.Cells.Item(18, 4).Value = 0.868
$ws.Cells.Item(19, 2).Value = 0.638
$ws.Cells.Item(20, 2).Value = 0.982
$ws.Cells.Item(20, 4).Value = 0.75
$ws.Cells.Item(20, 5).Value = 0.9330000000000001
$ws.Cells.Item(21, 2).Value = 0.52
$ws.Cells.Item(21, 4).Value = 0.781
$ws.Cells.Item(22, 2).Value = 1.112
$ws.Cells.Item(22, 3).Value = 0.004
$ws.Cells.Item(22, 4).Value = 0.626
$ws.Cells.Item(22, 5).Value = 0.533
$ws.Cells.Item(22, 6).Value = 'Andrea Ayunove Hutami'
$ws.Cells.Item(23, 2).Value = 1.075
$ws.Cells.Item(23, 3).Value = 0.004
$ws.Cells.Item(23, 4).Value = 0.795
$ws.Cells.Item(24, 2).Value = 1.154
$ws.Cells.Item(24, 3).Value = 0.004
$ws.Cells.Item(24, 4).Value = 0.882
$ws.Cells.Item(25, 2).Value = 1.213
$ws.Cells.Item(25, 3).Value = 0.004
$ws.Cells.Item(25, 4).Value = 0.906
$ws.Cells.Item(26, 2).Value = 0.9350000000000001
$ws.Cells.Item(26, 4).Value = 0.489
$ws.Cells.Item(26, 5).Value = 0.8
$ws.Cells.Item(27, 2).Value = 1.178
$ws.Cells.Item(27, 3).Value = 0.004
$ws.Cells.Item(27, 4).Value = 0.904
$ws.Cells.Item(28, 2).Value = 1.216
$ws.Cells.Item(28, 3).Value = 0.004
$ws.Cells.Item(28, 4).Value = 0.898
$ws.Cells.Item(29, 2).Value = 0.88
$ws.Cells.Item(29, 3).Value = 0.003
$ws.Cells.Item(30, 2).Value = 1.179
$ws.Cells.Item(30, 3).Value = 0.004
$ws.Cells.Item(30, 4).Value = 0.858
$ws.Cells.Item(31, 2).Value = 0.659
$ws.Cells.Item(31, 4).Value = 0.783
$ws.Cells.Item(31, 5).Value = 0.333
$ws.Cells.Item(32, 2).Value = 0.731
$ws.Cells.Item(32, 4).Value = 0.761
$ws.Cells.Item(33, 2).Value = 1.072
$ws.Cells.Item(33, 3).Value = 0.004
$ws.Cells.Item(33, 4).Value = 0.889
$ws.Cells.Item(33, 5).Value = 0.267
$ws.Cells.Item(33, 6).Value = 'Tidak Diketahui'
$ws.Cells.Item(33, 7).Value = 'Salah'
$ws.Cells.Item(34, 2).Value = 1.246
$ws.Cells.Item(34, 3).Value = 0.004
$ws.Cells.Item(34, 4).Value = 0.823
$ws.Cells.Item(34, 5).Value = 0.467
$ws.Cells.Item(34, 6).Value = 'Tidak Diketahui'
$ws.Cells.Item(34, 7).Value = 'Salah'
$ws.Cells.Item(35, 2).Value = 1.252
$ws.Cells.Item(35, 3).Value = 0.004
$ws.Cells.Item(35, 4).Value = 0.893
$ws.Cells.Item(35, 5).Value = 0.533
$ws.Cells.Item(36, 1).Value = 'FY_4.png'
$ws.Cells.Item(36, 2).Value = 1.159
$ws.Cells.Item(36, 3).Value = 0.004
$ws.Cells.Item(36, 4).Value = 0.8080000000000001
$ws.Cells.Item(36, 5).Value = 0.467
$ws.Cells.Item(36, 6).Value = 'Tidak Diketahui'
$ws.Cells.Item(36, 7).Value = 'Salah'
$ws.Cells.Item(37, 1).Value = 'TO_1.png'
$ws.Cells.Item(37, 2).Value = 0.79
$ws.Cells.Item(37, 4).Value = 0.792
$ws.Cells.Item(37, 5).Value = 0.8
$ws.Cells.Item(38, 1).Value = 'TO_2.png'
$ws.Cells.Item(38, 2).Value = 0.928
$ws.Cells.Item(38, 3).Value = 0.003
$ws.Cells.Item(38, 4).Value = 0.832
$ws.Cells.Item(39, 1).Value = 'TO_3.png'
$ws.Cells.Item(39, 2).Value = 0.8129999999999999
$ws.Cells.Item(39, 3).Value = 0.003
$ws.Cells.Item(39, 4).Value = 0.8149999999999999
$ws.Cells.Item(39, 5).Value = 1
$ws.Cells.Item(39, 6).Value = 'Tiara Oktavian'
$ws.Cells.Item(39, 7).Value = 'Benar'
$ws.Cells.Item(40, 1).Value = 'TO_4.png'
$ws.Cells.Item(40, 2).Value = 2.069
$ws.Cells.Item(40, 3).Value = 0.007
$ws.Cells.Item(40, 4).Value = 0.447
$ws.Cells.Item(40, 5).Value = 0.6
$ws.Cells.Item(41, 1).Value = 'TO_5.png'
$ws.Cells.Item(41, 2).Value = 1.952
$ws.Cells.Item(41, 3).Value = 0.006
$ws.Cells.Item(41, 4).Value = 0.463
$ws.Cells.Item(41, 7).Value = 'Salah'
$ws.Cells.Item(42, 1).Value = 'TD_1.png'
$ws.Cells.Item(42, 2).Value = 1.487
$ws.Cells.Item(42, 3).Value = 0.005
$ws.Cells.Item(42, 4).Value = 0.382
$ws.Cells.Item(43, 1).Value = 'TD_2.png'
$ws.Cells.Item(43, 2).Value = 1.474
$ws.Cells.Item(43, 3).Value = 0.005
$ws.Cells.Item(43, 4).Value = 0.362
$ws.Cells.Item(43, 5).Value = 0.4
$ws.Cells.Item(44, 1).Value = 'TD_3.png'
$ws.Cells.Item(44, 2).Value = 1.028
$ws.Cells.Item(44, 4).Value = 0.727
$ws.Cells.Item(44, 5).Value = 0.267
$ws.Cells.Item(45, 1).Value = 'TD_4.png'
$ws.Cells.Item(45, 2).Value = 0.969
$ws.Cells.Item(45, 3).Value = 0.003
$ws.Cells.Item(45, 4).Value = 0.695
$ws.Cells.Item(45, 5).Value = 0.267
$ws.Cells.Item(45, 6).Value = 'Tidak Diketahui'
$ws.Cells.Item(45, 7).Value = 'Benar'

Write-Output "Applied cell updates and appended row 45."
